# Append a new scraped listing (2025-09-28 18:22 JST) to the "ランサーズ"
# sheet: every existing row's "取得日時" timestamp is refreshed to the new
# run time, and a brand-new listing is inserted as the new row 4 (pushing
# the former rows 4-8 down to rows 5-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStamp = "2025-09-28 18:22:42"

# 1) Make room for the new listing right after the first two existing rows.
$ws.Rows.Item(4).Insert()

# 2) Populate the freshly inserted row 4 with the new listing's data.
$ws.Range("A4").Value = $newStamp
$ws.Range("B4").Value = "受付ソフトの自動入力システム開発依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5402230"
$ws.Range("G4").Value = 118
$ws.Range("H4").Value = "◆開発,システム開発"

# 3) Refresh the timestamp on every other data row (the insert already
#    shifted the old rows 4-8 down to rows 5-9).
$ws.Range("A2").Value = $newStamp
$ws.Range("A3").Value = $newStamp
$ws.Range("A5").Value = $newStamp
$ws.Range("A6").Value = $newStamp
$ws.Range("A7").Value = $newStamp
$ws.Range("A8").Value = $newStamp
$ws.Range("A9").Value = $newStamp

# 4) The row insert left the hyperlink relationships pointing at their old
#    (pre-shift) ranges, so rebuild the whole hyperlink collection for the
#    URL column against its final layout.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5394578")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5402230")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5402140")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5402038")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5402182")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5399347")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5402130")

# Adding hyperlinks through the COM layer stamps a brand-new style index on
# each cell; re-apply the shared "Hyperlink" cell style so the URL column
# keeps using the workbook's single existing hyperlink style.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
